$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Correct the 2021-05-26 (row 37) batch value, which recalculates the "% off" column
$ws.Range("B37").Value = 2223.4699999999998

# Append 2021-05-27 results: 5 samples from run 4 for TP2 (new row 38)
$ws.Range("A38").Value = 20210527
$ws.Range("B38").Value = 2221.3470000000002
$ws.Range("C38").Value = 2224.4699999999998
$ws.Range("D38").Formula = "=100*(B38-C38)/C38"
$ws.Range("E38").Value = 180
$ws.Range("F38").Value = "CRM opened 20210418"

# Leave the selection where the editor finished working
$ws.Range("E41").Select()
